$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "English" "Inglés"
Replace-Text " / Portuguese / French / Thai / Vietnamese / Spanish" " / Portugués / Francés / Tailandés / Vietnamita / Español"
Replace-Text "Brief" "Breve"
Replace-Text "An email sent upon verification to partners in the target country who have sent the correct documents. It will be sent via customer.io" "An email sent upon verification to partners in the target country who have sent the correct documents. Se enviará a través de customer.io"
Replace-Text "Target audience" "Público objetivo"
Replace-Text "Your documents have been verified!" "¡Tus documentos han sido verificados!"
Replace-Text "Hi " "Hola "
Replace-Text "We’ve reviewed the documents you’ve sent us for the " "Hemos revisado los documentos que nos has enviado para el "
Replace-Text " and all of them have been verified! " " y ¡todos han sido verificados! "
Replace-Text "We’ll be sending out more details about the event soon, including the agenda and travel arrangements, so make sure to check your inbox regularly." "Pronto te enviaremos más detalles sobre el evento, incluida la agenda y los preparativos del viaje, así que asegúrate de consultar regularmente tu bandeja de entrada."
Replace-Text "If you have any questions, please contact us via " "Si tienes alguna pregunta, entra en contacto con nosotros por "
Replace-Text "If you have any questions, please contact your country manager, " "Si tienes alguna pregunta, entra en contacto con el gestor de tu país "
Replace-Text ", at " ", en "
Replace-Text " or " " o "

foreach ($c in $d.Comments) {
    $c.Range.Find.Execute("choose either one", $true, $false, $false, $false, $false, `
                           $true, 1, $false, "elija uno de los dos", 2) | Out-Null
}
